$d = $word.ActiveDocument

# Locate the "{SomeCo, Inc.}," text and replace with "{SuperCo, Inc.},"
$d.Content.Find.Execute("{SomeCo, Inc.},", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{SuperCo, Inc.},", 2)

Write-Output "done"
